$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3: was stored as a number (2), becomes a text string "2"
$ws.Range("A3").Value = "'2"
$ws.Range("A3").Style = "Normal"

# E3: fix excessive spacing in the back of the card (bug fix per commit message)
$e3 = "<ul>`n<li>A etapa de desing</li>`n<li>são estabelecidos:`n<ul>`n`t<li>objetivos</li>`n`t<li>planos de ação</li>`n`t<li>bem como metas de qualidade</li>`n`t<li>os produtos e processos necessários à realização dessas metas</li>`n</ul>`n</li>`n</ul>`n"
$ws.Range("E3").Value = $e3

# The engine autofits row height on wrap-affecting edits; restore the
# original (unset) row height so only the cell contents change.
$ws.Rows(3).AutoFit()
